$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values must remain text, matching the original formatting
# (e.g. "9.10" with trailing zero, "61.576.34" with multiple separators).
# Force text number-format before assignment so Excel does not coerce these
# numeric-looking strings into real numbers, which would strip formatting.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.576.34'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.391.81'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.93'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.13'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.385'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.971.56'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.94'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.382.98'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.620.30'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.59'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.10'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '388.05'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.69'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.190'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.29'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.424.62'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.68'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.461.45'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.60'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0263'

# Column E (Volume/1h) values are percentage text with surrounding spaces;
# they are never numeric-parseable, so a plain text assignment is sufficient.
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('E14').Value = '  +1.08%  '
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('E20').Value = '  +2.01%  '
$ws.Range('E21').Value = '  +1.60%  '
$ws.Range('E23').Value = '  -0.92%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -3.83%  '
$ws.Range('E26').Value = '  +5.05%  '
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('E28').Value = '  +1.15%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  -0.73%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('E35').Value = '  +1.35%  '
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  +0.84%  '
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('E44').Value = '  +1.94%  '
$ws.Range('E46').Value = '  -2.04%  '
$ws.Range('E47').Value = '  -1.28%  '
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('E50').Value = '  -5.31%  '
$ws.Range('E51').Value = '  -1.05%  '
